# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The only substantive data change is that the "Salario Basico" values in
# column G for rows 21-25 (worker JUAN CARLOS DIAZ MERLANO, periods 2502-2506)
# drop from 1,430,000 to 1,423,500.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("G21:G25").Value = 1423500
